$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new item row populated
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "ANTI-COX II 15MG/3ML 6 AMP"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "3:4"
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "78.00"
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "13.2600"
$ws.Range("Q7").Value = "0:1"

# Row 8: totals row now has a value + taller row
$ws.Range("N8").Value = 13.26
$ws.Rows("8").RowHeight = 25.5

# Row 9: refreshed footer (timestamp / page / developer credit)
$ws.Range("A9").Value = "Monday, 29 December, 2025 9:05 AM"
$ws.Range("G9").Value = "1/1"
$ws.Range("K9").Value = "developed by : Abdelaziz Talaat"
